# The workbook stores every field as text, even values that look numeric
# (phone numbers, years, dates as plain strings, etc). Row 25's "Phone"
# cell (F25) had accidentally been entered as a formula (=8420880979)
# instead of the plain text "08420880979" that every other row uses.
# This change removes that stray formula, and appends a new registrant
# (row 26) - a duplicate "AGNIVA BHATTACHARJEE" entry with a changed
# email/content and its own Phone cell re-using the same formula-style
# entry pattern the sheet already had in row 25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix F25: drop the stray "=8420880979" formula, keep it as the
#     plain text "08420880979" (NumberFormat "@" keeps the leading zero
#     and stops it being re-interpreted as a number). ---
$ws.Range("F25").ClearContents()
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = "08420880979"

# --- Append new row 26 ---
$ws.Range("A26").Value = "BSS/d32409a3c007"
$ws.Range("B26").Value = "AGNIVA"
$ws.Range("C26").Value = "BHATTACHARJEE"
$ws.Range("D26").Value = "kuch nhi hai"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2020"

# Phone: same "formula that evaluates to the phone number" pattern the
# original sheet used for this field.
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Formula = "=8420880979"

$ws.Range("G26").Value = "bhattacharjee.agniva.jobs@gmail.com"

$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "2002-01-21"

$ws.Range("I26").Value = "IT"
$ws.Range("J26").Value = "Google"

$ws.Range("K26").Value = ""
$ws.Range("L26").Value = ""
$ws.Range("M26").Value = ""

$ws.Range("N26").Value = "456465kjhgfg"
